$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these numeric-looking cells remain stored as text, matching the source data (inline strings)
$textCells = @("D5", "D6", "D10", "D14", "D18", "D20", "D21", "D22", "D27", "D28", "D30", "D31", "D33", "D35", "D36", "D38", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.794.61"
$ws.Range("E2").Value = "  +3.00%  "
$ws.Range("D3").Value = "3.384.36"
$ws.Range("E3").Value = "  +3.81%  "
$ws.Range("D5").Value = "191.64"
$ws.Range("E5").Value = "  +3.43%  "
$ws.Range("D6").Value = "594.22"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("D10").Value = "6.76"
$ws.Range("E10").Value = "  +2.65%  "
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("D12").Value = "3.975.19"
$ws.Range("E12").Value = "  +3.98%  "
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").Value = "28.71"
$ws.Range("E14").Value = "  +3.48%  "
$ws.Range("D15").Value = "69.761.67"
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").Value = "3.377.24"
$ws.Range("E17").Value = "  +2.26%  "
$ws.Range("D18").Value = "454.63"
$ws.Range("E18").Value = "  +15.36%  "
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").Value = "13.84"
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("D21").Value = "7.80"
$ws.Range("E21").Value = "  +2.38%  "
$ws.Range("D22").Value = "76.08"
$ws.Range("E22").Value = "  +6.32%  "
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("E25").Value = "  +3.27%  "
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("D27").Value = "9.54"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  +3.57%  "
$ws.Range("D30").Value = "23.46"
$ws.Range("E30").Value = "  +3.36%  "
$ws.Range("D31").Value = "5.61"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").Value = "7.01"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +6.69%  "
$ws.Range("D36").Value = "164.68"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("E37").Value = "  +2.30%  "
$ws.Range("D38").Value = "27.93"
$ws.Range("E38").Value = "  +4.69%  "
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").Value = "4.60"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").Value = "6.61"
$ws.Range("E41").Value = "  +2.01%  "
$ws.Range("D42").Value = "2.745.96"
$ws.Range("E42").Value = "  +5.01%  "
$ws.Range("D43").Value = "2.52"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("D44").Value = "25.54"
$ws.Range("E44").Value = "  +2.58%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "41.16"
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0688"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("D47").Value = "339.82"
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("D48").Value = "0.0284"
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("D49").Value = "32.90"
$ws.Range("E49").Value = "  +6.58%  "
$ws.Range("E50").Value = "  +4.54%  "
$ws.Range("E51").Value = "  -0.35%  "
